# Adds tunable parameters for control:
#  - "Mass" sheet gains a motors/frame weight breakdown (rows 12-13)
#  - "Test" chart's polynomial trendline gains a companion linear trendline
#  - View-state (selection / zoom / scroll) for both sheets is updated

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Mass sheet: new "motors" / "frame" rows feeding the Total (D2) sum
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("Mass")

$ws2.Range("A12").Value = "motors"
$ws2.Range("B12").Value = 36
$ws2.Range("C12").Formula = "=B12*4"
$ws2.Range("D12").Formula = "=C12+B13"

$ws2.Range("A13").Value = "frame"
$ws2.Range("B13").Value = 132

# ---------------------------------------------------------------------
# Test sheet chart: add a second (linear) trendline alongside the
# existing polynomial one, with R^2 and the fitted equation shown.
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("Test")
$chartObj = $ws1.ChartObjects().Item(1)
$chart = $chartObj.Chart
$series = $chart.SeriesCollection().Item(1)
$trendlines = $series.Trendlines()
# xlLinear = -4132; DisplayEquation/DisplayRSquared both on, as in the workbook
$trendlines.Add(-4132, 2, 0, 0, 0, $true, $true, "Linear (Test)")
$newTrend = $trendlines.Item($trendlines.Count)
$newTrend.DisplayEquation = $true
$newTrend.DisplayRSquared = $true

# ---------------------------------------------------------------------
# View state: selection / zoom / scroll per sheet. Touch "Mass" first so
# the workbook ends up with "Test" as the active sheet/tab, matching the
# saved file.
# ---------------------------------------------------------------------
[void]$ws2.Activate()
$excel.ActiveWindow.Zoom = 181
[void]$ws2.Range("D15").Select()

[void]$ws1.Activate()
$excel.ActiveWindow.ScrollColumn = 4
$excel.ActiveWindow.ScrollRow = 1
[void]$ws1.Range("U11").Select()
